$d = $word.ActiveDocument

# The count-API paragraph originally read "...and call the helper method
# getUrl and if the URL exists..." -- update it to reference the new
# repository method used by the delete/count logic.
$d.Content.Find.Execute("helper method getUrl", $true, $false, $false, $false, $false,
                         $true, 1, $false, "findByUrl in the repository", 2)
